# Update the Rules worksheet: change the "Good Morning" greeting cell (E8)
# to "GIT UPDATE", matching the new shared-string entry, and leave the
# selection on E8 (as recorded in the saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
